$d = $word.ActiveDocument

$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)

# The document currently ends with a single empty paragraph (just a
# paragraph mark, right before the sectPr). Push a fresh paragraph mark
# in after it, then delete the *original* empty paragraph's own mark so
# it disappears, leaving only the brand-new mark behind as the (now)
# last paragraph of the body. This keeps the preceding "..." paragraph's
# identity completely untouched.
$rng = $lastPara.Range
$rng.Collapse(1)
$rng.InsertParagraphAfter()

$origEmpty = $d.Paragraphs.Item($paraCount)
$markRng = $d.Range($origEmpty.Range.End - 1, $origEmpty.Range.End)
$markRng.Delete()

# Insert a second brand-new paragraph at the very end of the document.
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

$newCount = $d.Paragraphs.Count
$firstNew = $d.Paragraphs.Item($newCount - 1)
$secondNew = $d.Paragraphs.Item($newCount)

$firstNew.Range.Text = 'Version management is required during project development or updating the version of a developed project. Normally, when a group of people work on same project, there could be clashes on development as different changes are made by different users. To manage these changes on a single project, a version management software like Git and GitHub is used by project managers. GitHub is repositories for the projects where version is managed by creating various branches. The “Main” is the mainline of the project in the repositories and the finalized component of software or project is stored in the main branch while, many other sub-branches could be assigned to each developer. Each developer will develop assigned component of the project and push to the assigned branch. Once the branch is fully tested and approved, then the project manager will merge it to main branch. '

$secondNew.Range.Text = 'For functional projects, updated version are pushed to branch while main line will consist working or functional piece of project or software. Once the updated version is fully tested and approved then, it is merged to the main-line as the new version of the software.'
